$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7, shifting the existing rows 7-9 down to 8-10.
$ws.Rows.Item(7).Insert()

# Fill the newly inserted row 7 with the new Automobile Insurance field-hint/error test entry.
$ws.Range("A7").Value = "102_AutomobileInsurance_003_InsurantData_002_FieldHintsAndErrors"
$ws.Range("B7").Value = "var102_AutomobileInsurance_003_InsurantData_002_FieldHintsAndErrors"
$ws.Range("C7").Value = "Open Automobile Insurance"
$ws.Range("D7").Value = "102_AutomobileInsurance_003_InsurantData_002_FieldHintsAndErrors"

# Update the active selection, as left by the author's Excel session.
$ws.Range("A7").Select()
